# Apply updates to row 4 (match: Ind. Medellin vs Llaneros) odds values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 1.27
$ws.Range("H4").Value = 5
$ws.Range("J4").Value = 1.73
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 10
$ws.Range("Q4").Value = 1.85
$ws.Range("R4").Value = 2
$ws.Range("S4").Value = 3
$ws.Range("T4").Value = 1.36
$ws.Range("U4").Value = 1.36
$ws.Range("V4").Value = 3
$ws.Range("Y4").Value = 6
$ws.Range("AE4").Value = 10
$ws.Range("AF4").Value = 10
$ws.Range("AR4").Value = 2.48
